# Update "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" sheets, as published in the regenerated
# gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row number -> new F-column value }
$updates = @{
    "展览" = @{
        4  = 3737
        5  = 2265
        6  = 443
        9  = 182
        10 = 102
        12 = 1383
        14 = 2206
        15 = 158
    }
    "全部类型" = @{
        4  = 3737
        5  = 2265
        6  = 443
        10 = 182
        11 = 102
        15 = 1383
        17 = 2206
        18 = 158
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
